$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REIT")

# Fill in the previously-blank monthly NAV / distribution rows (Apr 2020 - Nov 2020)
$ws.Range("A131").Value = 43922
$ws.Range("B131").Value = 18.896000000000001
$ws.Range("C131").Value = 0.068330000000000002

$ws.Range("A132").Value = 43952
$ws.Range("B132").Value = 18.896000000000001
$ws.Range("C132").Value = 0.068330000000000002

$ws.Range("A133").Value = 43983
$ws.Range("B133").Value = 18.896000000000001
$ws.Range("C133").Value = 0.068330000000000002

$ws.Range("A134").Value = 44013
$ws.Range("B134").Value = 18.896000000000001
$ws.Range("C134").Value = 0.068330000000000002

$ws.Range("A135").Value = 44044
$ws.Range("B135").Value = 19.242999999999999
$ws.Range("C135").Value = 0.068330000000000002

$ws.Range("A136").Value = 44075
$ws.Range("B136").Value = 19.242999999999999
$ws.Range("C136").Value = 0.068330000000000002

$ws.Range("A137").Value = 44105
$ws.Range("B137").Value = 19.242999999999999
$ws.Range("C137").Value = 0.068330000000000002

$ws.Range("A138").Value = 44136
$ws.Range("B138").Value = 19.34
$ws.Range("C138").Value = 0.068330000000000002

# Match the saved cursor / selection state
[void]$ws.Range("A124").Select()
[void]$ws.Range("C141").Select()
